# Append the new match row (row 61) to the betexplorer sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 61
$prevRow = 60

# Copy the formatting of the previous data row down to the new row first,
# so number formats / borders / font for the "Indice" and date columns match.
$srcRange = $ws.Range("A" + $prevRow + ":V" + $prevRow)
$dstRange = $ws.Range("A" + $newRow + ":V" + $newRow)
$srcRange.Copy()
$dstRange.PasteSpecial(-4122)  # xlPasteFormats

$rowData = @(
    60,
    "algeria",
    "ligue-1",
    "2023-2024",
    45267.79166666666,
    "Paradou",
    0,
    "MC Alger",
    1,
    3.44,
    "06/12/2023 07:12",
    8.82,
    "07/12/2023 18:38",
    2.7,
    "06/12/2023 07:12",
    3.82,
    "07/12/2023 18:38",
    2.3,
    "06/12/2023 07:12",
    1.47,
    "07/12/2023 18:38",
    "https://www.betexplorer.com/football/algeria/ligue-1/paradou-mc-alger/ITA1dIkm/"
)

for ($col = 1; $col -le $rowData.Length; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $rowData[$col - 1]
}
